$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''62.332.80'
$ws.Range("E2").Value = '  +0.65%  '
$ws.Range("D3").Value = '''3.430.81'
$ws.Range("E3").Value = '  +0.42%  '
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("D5").Value = '''413.46'
$ws.Range("E5").Value = '  +0.87%  '
$ws.Range("D6").Value = '''128.74'
$ws.Range("E6").Value = '  -0.25%  '
$ws.Range("D7").Value = '''0.625'
$ws.Range("E7").Value = '  -2.23%  '
$ws.Range("E8").Value = '  +0.05%  '
$ws.Range("D9").Value = '''0.726'
$ws.Range("E9").Value = '  -1.30%  '
$ws.Range("E10").Value = '  +0.50%  '
$ws.Range("E11").Value = '  -0.30%  '
$ws.Range("D12").Value = '''0.0000221'
$ws.Range("E12").Value = '  +9.05%  '
$ws.Range("D13").Value = '''9.22'
$ws.Range("E13").Value = '  +1.16%  '
$ws.Range("D14").Value = '''3.971.55'
$ws.Range("E14").Value = '  +0.71%  '
$ws.Range("E15").Value = '  -0.27%  '
$ws.Range("D16").Value = '''20.45'
$ws.Range("E16").Value = '  -3.47%  '
$ws.Range("D17").Value = '''3.413.17'
$ws.Range("E17").Value = '  +0.70%  '
$ws.Range("D18").Value = '''12.73'
$ws.Range("E18").Value = '  +5.00%  '
$ws.Range("E19").Value = '  -1.02%  '
$ws.Range("D20").Value = '''62.334.92'
$ws.Range("E20").Value = '  +0.95%  '
$ws.Range("D21").Value = '''476.44'
$ws.Range("E21").Value = '  +4.34%  '
$ws.Range("D22").Value = '''91.76'
$ws.Range("E22").Value = '  -0.46%  '
$ws.Range("E23").Value = '  +3.07%  '
$ws.Range("D24").Value = '''13.12'
$ws.Range("E24").Value = '  +1.48%  '
$ws.Range("D25").Value = '''3.32'
$ws.Range("E25").Value = '  +2.08%  '
$ws.Range("D26").Value = '''9.69'
$ws.Range("E26").Value = '  +10.23%  '
$ws.Range("D27").Value = '''33.41'
$ws.Range("E27").Value = '  -1.18%  '
$ws.Range("D28").Value = '''4.78'
$ws.Range("E28").Value = '  +0.89%  '
$ws.Range("D29").Value = '''7.71'
$ws.Range("E29").Value = '  +1.47%  '
$ws.Range("D30").Value = '''11.88'
$ws.Range("E30").Value = '  -1.08%  '
$ws.Range("D31").Value = '''2.63'
$ws.Range("E31").Value = '  -4.36%  '
$ws.Range("D32").Value = '''0.166'
$ws.Range("E32").Value = '  -1.14%  '
$ws.Range("E33").Value = '  -2.39%  '
$ws.Range("D34").Value = '''40.78'
$ws.Range("E34").Value = '  -4.83%  '
$ws.Range("D36").Value = '''58.02'
$ws.Range("E36").Value = '  +8.68%  '
$ws.Range("E37").Value = '  -2.00%  '
$ws.Range("D38").Value = '''0.999'
$ws.Range("E38").Value = '  +0.28%  '
$ws.Range("D39").Value = '''3.02'
$ws.Range("E39").Value = '  +3.78%  '
$ws.Range("E40").Value = '  +0.39%  '
$ws.Range("E41").Value = '  +2.68%  '
$ws.Range("D42").Value = '''3.33'
$ws.Range("E42").Value = '  -1.35%  '
$ws.Range("D43").Value = '''2.65'
$ws.Range("E43").Value = '  +10.09%  '
$ws.Range("D44").Value = '''143.94'
$ws.Range("E44").Value = '  +2.17%  '
$ws.Range("B45").Value = 'ARBITRUM'
$ws.Range("C45").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D45").Value = '''2.05'
$ws.Range("E45").Value = '  +3.83%  '
$ws.Range("B46").Value = 'NEARProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D46").Value = '''4.30'
$ws.Range("E46").Value = '  +2.51%  '
$ws.Range("D47").Value = '''2.41'
$ws.Range("E47").Value = '  +19.74%  '
$ws.Range("D48").Value = '''16.32'
$ws.Range("E48").Value = '  -1.68%  '
$ws.Range("E49").Value = '  +31.17%  '
$ws.Range("D50").Value = '''22.41'
$ws.Range("E50").Value = '  -0.18%  '
$ws.Range("D51").Value = '''112.58'
$ws.Range("E51").Value = '  +5.17%  '
